$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Φύλλο1" to "1"
$ws.Name = "1"

# Feature engineering from timestamp: the first table column ("Nama Gedung")
# is repurposed to hold the meter_id instead of the building name, so the
# header is renamed accordingly (the table column name follows automatically).
$ws.Range("A1").Value = "meter_id"

# Trim the trailing space in the "SBM " building-name entry.
$ws.Range("A39").Value = "SBM"

# Correct the is_kelas / is_penelitian flags for "Gedung Laboratorium Doping".
$ws.Range("B13").Value = 0
$ws.Range("D13").Value = 1

# Restore the scroll position / active selection of the sheet view.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("I39").Select()
